# Shift each data row (2..20) one column to the left (B<-C, C<-D, ...),
# dropping the old value that was in column B. Rows 2..10 pick up a brand
# new trailing data point (the newly available ifoCAST observation for that
# quarter) in the now-vacated last column; rows 11..20 simply lose their
# last column with nothing to replace it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data points appended at the end of rows 2-10 (column K) - the newly
# evaluated ifoCAST quarter that extends each row's error series by one.
$newLastValues = @{
    2  = 0.3406326034063205
    3  = -1.53061581027015
    4  = 1.738194353941456
    5  = 1.952515937864398
    6  = -0.3663687737149753
    7  = 0.05323806995971928
    8  = -0.1665195462441563
    9  = 0.6100525277605273
    10 = -0.4353061035472806
}

for ($row = 2; $row -le 20; $row++) {
    # Find the last used column in this row (data starts at column 2 = B).
    $lastCol = 1
    for ($col = 2; $col -le 11; $col++) {
        $val = $ws.Cells.Item($row, $col).Value2
        if ($val -ne $null -and $val -ne "") {
            $lastCol = $col
        }
    }

    if ($lastCol -ge 3) {
        # Read the old values first (columns 2..lastCol) ...
        $oldVals = @{}
        for ($col = 2; $col -le $lastCol; $col++) {
            $oldVals[$col] = $ws.Cells.Item($row, $col).Value2
        }
        # ... then write them back shifted one column to the left.
        for ($col = 2; $col -le ($lastCol - 1); $col++) {
            $ws.Cells.Item($row, $col).Value = $oldVals[$col + 1]
        }
    }

    if ($newLastValues.ContainsKey($row)) {
        # This row keeps its full width - fill the newly freed last cell
        # with the new trailing data point.
        $ws.Cells.Item($row, $lastCol).Value = $newLastValues[$row]
    } else {
        # This row shrinks by one - clear the now-unused last cell.
        $ws.Cells.Item($row, $lastCol).ClearContents()
    }
}
